$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S50003.MES.BIN")

# --- Edited-translation text fixes (in-place replacements of existing strings) ---
$ws.Range("G2").Value = "Huh? What is?"
$ws.Range("G3").Value = "Without warning, I ask her without thinking."
$ws.Range("G6").Value = "See you later, Fujii-san. Thanks for your hard work today. Bye bye."

# --- Fill in previously-empty "Edited" column cells (column F) ---
$ws.Range("F7").Value = "Without saying anything, I follow her out of the house, carried away by her flowing momentum."
$ws.Range("F8").Value = "Well, I'll be going. Goodbye, Fujii-san."

# --- Replace Edited column text for row 9 ---
$ws.Range("F9").Value = "She was so sullen when she was facing the desk, but she's a girl with rapid mood changes."

# --- Reset the lingering cell selection (F8) back to the top-left cell ---
$ws.Range("A1").Select()
